# Applies the daily "remaining days" refresh to Sheet1.
# For every data row (2-95) except row 36 (which has a malformed date
# and was left untouched in the source change), column E ("剩余" / days
# remaining) is decremented by 1, reflecting one more day having elapsed.
# Row 95 additionally has its start date (column F) rolled forward to
# 20251012, which resets its remaining-days count to 10 (its full total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skipRows = @(36)

for ($r = 2; $r -le 95; $r++) {
    if ($skipRows -contains $r) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Value = $cell.Value() - 1
}

# Row 95: the reference start date moved to 2025-10-12, resetting its
# remaining-day count back up to its total (column D = 10).
$ws.Cells.Item(95, 6).Value = 20251012   # column F
$ws.Cells.Item(95, 5).Value = 10         # column E
